# Updates the "Estado de Cuenta" data table (rows 16-35) on sheet Hoja1:
#  - Column E (Periodo Mora) values are refreshed with a new set of periods
#  - Column F (Valor Mora) values are refreshed
#  - Column G (Salario Basico) is updated to the new base salary value
#
# This mirrors the macro-driven "se elimina EC anteriores y se agregan
# nuevos, se modifica base de datos" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row=16; Periodo="2009"; Mora=32771; Salario=877803 },
    @{ Row=17; Periodo="2008"; Mora=35112; Salario=877803 },
    @{ Row=18; Periodo="2007"; Mora=35112; Salario=877803 },
    @{ Row=19; Periodo="2006"; Mora=35112; Salario=877803 },
    @{ Row=20; Periodo="2005"; Mora=35112; Salario=877803 },
    @{ Row=21; Periodo="2004"; Mora=35112; Salario=877803 },
    @{ Row=22; Periodo="2002"; Mora=33125; Salario=877803 },
    @{ Row=23; Periodo="2001"; Mora=33125; Salario=877803 },
    @{ Row=24; Periodo="1911"; Mora=33125; Salario=877803 },
    @{ Row=25; Periodo="1910"; Mora=33125; Salario=877803 },
    @{ Row=26; Periodo="1909"; Mora=33125; Salario=877803 },
    @{ Row=27; Periodo="1907"; Mora=31249; Salario=877803 },
    @{ Row=28; Periodo="1906"; Mora=31249; Salario=877803 },
    @{ Row=29; Periodo="1905"; Mora=31249; Salario=877803 },
    @{ Row=30; Periodo="1904"; Mora=31249; Salario=877803 },
    @{ Row=31; Periodo="1903"; Mora=31249; Salario=877803 },
    @{ Row=32; Periodo="1902"; Mora=31249; Salario=877803 },
    @{ Row=33; Periodo="1808"; Mora=31249; Salario=877803 },
    @{ Row=34; Periodo="1807"; Mora=31249; Salario=877803 },
    @{ Row=35; Periodo="1805"; Mora=31249; Salario=877803 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
